$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.830.15"
$ws.Range("E2").Value = "  +3.44%  "
$ws.Range("D3").Value = "3.413.54"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.91"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.83"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("D8").Value = "3.406.83"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +13.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.633"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.19"
$ws.Range("E12").Value = "  +3.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000281"
$ws.Range("E13").Value = "  +5.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.17"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "3.958.88"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.35"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.119"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.406.42"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "65.771.38"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.90"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.97"
$ws.Range("E22").Value = "  +14.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"
$ws.Range("E23").Value = "  +20.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.16"
$ws.Range("E24").Value = "  +2.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.68"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.51"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  +3.56%  "
$ws.Range("E28").Value = "  +6.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.91"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.01"
$ws.Range("E30").Value = "  +6.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.73"
$ws.Range("E31").Value = "  +4.99%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.57"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.86"
$ws.Range("E33").Value = "  +9.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "579.78"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  +4.27%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.87"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("D42").Value = "3.100.47"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.136"
$ws.Range("E47").Value = "  +6.19%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.22"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.39"
$ws.Range("E50").Value = "  +5.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.64"
$ws.Range("E51").Value = "  +2.67%  "
